$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 54, pushing the existing rows 54:170 down to 56:172
$ws.Rows("54:55").Insert()

# Populate the newly inserted row 54 with the new week's "Especial" quality record
$ws.Range("A54").Value = 4
$ws.Range("B54").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value = "Los Lagos"
$ws.Range("D54").Value = 44526
$ws.Range("E54").Value = 10
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100101
$ws.Range("H54").Value = "Berries"
$ws.Range("I54").Value = 100101007
$ws.Range("J54").Value = "Kiwi"
$ws.Range("K54").Value = "Hayward"
$ws.Range("L54").Value = "Especial"
$ws.Range("M54").Value = 200
$ws.Range("N54").Value = 22000
$ws.Range("O54").Value = 22000
$ws.Range("P54").Value = 22000
$ws.Range("Q54").Value = "`$/caja 15 kilos"
$ws.Range("R54").Value = "Provincia de Curicó"
$ws.Range("S54").Value = 1467
$ws.Range("T54").Value = 15

# Populate the newly inserted row 55 with the new week's "Primera" quality record
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C55").Value = "Los Lagos"
$ws.Range("D55").Value = 44526
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100101
$ws.Range("H55").Value = "Berries"
$ws.Range("I55").Value = 100101007
$ws.Range("J55").Value = "Kiwi"
$ws.Range("K55").Value = "Hayward"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 400
$ws.Range("N55").Value = 16000
$ws.Range("O55").Value = 16000
$ws.Range("P55").Value = 16000
$ws.Range("Q55").Value = "`$/caja 15 kilos"
$ws.Range("R55").Value = "Provincia de Curicó"
$ws.Range("S55").Value = 1067
$ws.Range("T55").Value = 15

# Keep the date column's custom format consistent with the rest of the column
$ws.Range("D54:D55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
